$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# Updated "First get config" RPC reply (F2) - now contains full network-instance data
$f2 = "<rpc-reply message-id=`"urn:uuid:9321e9f2-64b4-4653-865d-5a419ed77ba7`">`n"
$f2 += "  <data>`n"
$f2 += "    <network-instances>`n"
$f2 += "      <network-instance>`n"
$f2 += "        <name>Prueba_LxVPN</name>`n"
$f2 += "        <config>`n"
$f2 += "          <name>Prueba_LxVPN</name>`n"
$f2 += "          <type>oc-ni-types:L3VRF</type>`n"
$f2 += "        </config>`n"
$f2 += "        <interfaces>`n"
$f2 += "          <interface>`n"
$f2 += "            <id>GigabitEthernet0/3/2</id>`n"
$f2 += "            <config>`n"
$f2 += "              <id>GigabitEthernet0/3/2</id>`n"
$f2 += "              <interface>GigabitEthernet0/3/2</interface>`n"
$f2 += "              <subinterface>0</subinterface>`n"
$f2 += "            </config>`n"
$f2 += "          </interface>`n"
$f2 += "        </interfaces>`n"
$f2 += "        <protocols>`n"
$f2 += "          <protocol>`n"
$f2 += "            <identifier>oc-pol-types:OSPF</identifier>`n"
$f2 += "            <name>22</name>`n"
$f2 += "            <config>`n"
$f2 += "              <identifier>oc-pol-types:OSPF</identifier>`n"
$f2 += "              <name>22</name>`n"
$f2 += "            </config>`n"
$f2 += "            <ospfv2>`n"
$f2 += "              <global>`n"
$f2 += "                <config>`n"
$f2 += "                  <router-id>172.16.1.3</router-id>`n"
$f2 += "                </config>`n"
$f2 += "              </global>`n"
$f2 += "            </ospfv2>`n"
$f2 += "          </protocol>`n"
$f2 += "          <protocol>`n"
$f2 += "            <identifier>oc-pol-types:STATIC</identifier>`n"
$f2 += "            <name>default</name>`n"
$f2 += "            <config>`n"
$f2 += "              <identifier>oc-pol-types:STATIC</identifier>`n"
$f2 += "              <name>default</name>`n"
$f2 += "            </config>`n"
$f2 += "          </protocol>`n"
$f2 += "          <protocol>`n"
$f2 += "            <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>`n"
$f2 += "            <name>default</name>`n"
$f2 += "            <config>`n"
$f2 += "              <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>`n"
$f2 += "              <name>default</name>`n"
$f2 += "            </config>`n"
$f2 += "          </protocol>`n"
$f2 += "        </protocols>`n"
$f2 += "      </network-instance>`n"
$f2 += "    </network-instances>`n"
$f2 += "  </data>`n"
$f2 += "</rpc-reply>`n"

$ws.Range("F2").Value = $f2

# Updated "RPC" edit-config request (G2) - protocol name/identifier and interface id changed
$g2 = "<edit-config>`n"
$g2 += "    <target>`n"
$g2 += "      <candidate/>`n"
$g2 += "    </target>`n"
$g2 += "    <config>`n"
$g2 += "      <network-instances xmlns=`"http://openconfig.net/yang/network-instance`">`n"
$g2 += "        <network-instance>`n"
$g2 += "          <name>Prueba_LxVPN</name>`n"
$g2 += "          <config>`n"
$g2 += "            <name>Prueba_LxVPN</name>`n"
$g2 += "            <type xmlns:oc-ni-types=`"http://openconfig.net/yang/network-instance-types`">oc-ni-types:L3VRF</type>`n"
$g2 += "          </config>`n"
$g2 += "          <protocols>`n"
$g2 += "            <protocol>`n"
$g2 += "              <identifier xmlns:oc-pol-types=`"http://openconfig.net/yang/policy-types`">oc-pol-types:OSPF</identifier>`n"
$g2 += "              <name>22</name>`n"
$g2 += "              <config>`n"
$g2 += "                <identifier xmlns:oc-pol-types=`"http://openconfig.net/yang/policy-types`">oc-pol-types:OSPF</identifier>`n"
$g2 += "                <name>22</name>`n"
$g2 += "              </config>`n"
$g2 += "              <ospfv2>`n"
$g2 += "                <areas>`n"
$g2 += "                  <area>`n"
$g2 += "                    <identifier>0.0.0.0</identifier>`n"
$g2 += "                    <config>`n"
$g2 += "                      <identifier>0.0.0.0</identifier>`n"
$g2 += "                    </config>`n"
$g2 += "                    <interfaces>`n"
$g2 += "                      <interface>`n"
$g2 += "                        <id>GigabitEthernet0/3/0</id>`n"
$g2 += "                        <config>`n"
$g2 += "                          <id>GigabitEthernet0/3/0</id>`n"
$g2 += "                        </config>`n"
$g2 += "                      </interface>`n"
$g2 += "                    </interfaces>`n"
$g2 += "                  </area>`n"
$g2 += "                </areas>`n"
$g2 += "              </ospfv2>`n"
$g2 += "            </protocol>`n"
$g2 += "          </protocols>`n"
$g2 += "        </network-instance>`n"
$g2 += "      </network-instances>`n"
$g2 += "    </config>`n"
$g2 += "</edit-config>"

$ws.Range("G2").Value = $g2
